$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells hold plain text in the source data (prices/percentages
# formatted as strings, some with multiple dots or leading zeros), so force
# text format before assigning to prevent Excel from auto-converting them
# to numbers and losing the exact original formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.789.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.465.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.25"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.17"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.156"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.84"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "68.694.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000170"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.53"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.61"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "336.77"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.88"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.58"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.65"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.23"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0₃0820"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.20"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "427.50"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.14"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.62"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.96"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.82"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.43"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.07"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.05"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.38"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "130.10"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0720"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.38"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.99"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -7.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.82"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0205"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.32%  "
